$wb = $excel.ActiveWorkbook

# --- hotel_info: insert a new "State" column between "Hotel_Name" and "City" ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- reorder tabs: review_info becomes the first sheet, hotel_info the second ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($wb.Worksheets.Item(1))
